# Actualización automática de noticias - 2026-01-18
# Insert a new row at the top of the data (row 2), pushing the existing
# rows down by one, and populate it with the newest news item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts old rows 2-15 down to 3-16)
$ws.Rows.Item(2).Insert()

# The inserted row inherits the formatting of the row above (header row);
# clear it so the new row looks like a plain data row, same as the others.
$ws.Range("A2:F2").ClearFormats()

# Leading apostrophe keeps the date-like text as a plain string instead of
# letting Excel auto-convert it into a date value (matches the other rows,
# which all store "fecha" as inline text, e.g. "2026-01-17").
$ws.Range("A2").Value = "'2026-01-18"
$ws.Range("B2").Value = "Este es el calendario de vacaciones y recesos en colegios públicos de Bogotá en 2026, prográmese"
$ws.Range("C2").Value = "Infobae"
$ws.Range("D2").Value = "Bogotá"
$ws.Range("E2").Value = "https://www.infobae.com/colombia/2026/01/17/este-es-el-calendario-de-vacaciones-y-recesos-en-colegios-publicos-de-bogota-en-2026-programese/"
$ws.Range("F2").Value = ""
